$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (price) holds text-formatted price strings (e.g. '28.277.77',
# '0.9985') that must stay text, not get auto-converted to numbers by Excel's
# input parsing. Restrict to the data rows so the header cell D1 keeps its
# original style untouched.
$ws.Range('D2:D51').NumberFormat = '@'

# Row 2: Bitcoin
$ws.Range('D2').Value = '28.277.77'
$ws.Range('E2').Value = '  -0.74%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.808.73'
$ws.Range('E3').Value = '  -0.91%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.63%  '

# Row 5: BNB
$ws.Range('D5').Value = '313.89'
$ws.Range('E5').Value = '  -0.81%  '

# Row 6: USDC
$ws.Range('D6').Value = '0.9981'
$ws.Range('E6').Value = '  -0.57%  '

# Row 7: XRP
$ws.Range('D7').Value = '0.5160'
$ws.Range('E7').Value = '  -0.16%  '

# Row 8: Cardano
$ws.Range('D8').Value = '0.3983'
$ws.Range('E8').Value = '  +3.07%  '

# Row 9: Dogecoin
$ws.Range('D9').Value = '0.07879'
$ws.Range('E9').Value = '  -5.53%  '

# Row 10: Polygon
$ws.Range('E10').Value = '  -0.72%  '

# Row 11: OKB
$ws.Range('D11').Value = '41.12'
$ws.Range('E11').Value = '  -2.02%  '

# Row 12: Polkadot
$ws.Range('D12').Value = '6.344'
$ws.Range('E12').Value = '  -1.08%  '

# Row 13: BinanceUSD
$ws.Range('D13').Value = '0.9979'
$ws.Range('E13').Value = '  -0.64%  '

# Row 14: Solana
$ws.Range('D14').Value = '20.43'
$ws.Range('E14').Value = '  -3.52%  '

# Row 15: Chainlink
$ws.Range('D15').Value = '7.329'
$ws.Range('E15').Value = '  -2.25%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '1.798.48'
$ws.Range('E16').Value = '  -1.36%  '

# Row 17: Litecoin
$ws.Range('D17').Value = '92.79'
$ws.Range('E17').Value = '  -1.37%  '

# Row 18: ShibaInu
$ws.Range('E18').Value = '  -3.85%  '

# Row 19: TRON
$ws.Range('D19').Value = '0.06571'
$ws.Range('E19').Value = '  -1.02%  '

# Row 20: Dai
$ws.Range('D20').Value = '0.9975'
$ws.Range('E20').Value = '  -0.59%  '

# Row 21: Avalanche
$ws.Range('E21').Value = '  -2.35%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '6.010'
$ws.Range('E22').Value = '  -0.94%  '

# Row 23: WrappedBTC
$ws.Range('D23').Value = '28.318.72'
$ws.Range('E23').Value = '  -0.77%  '

# Row 24: Cosmos
$ws.Range('D24').Value = '11.15'
$ws.Range('E24').Value = '  -2.18%  '

# Row 25: Toncoin
$ws.Range('D25').Value = '2.224'
$ws.Range('E25').Value = '  -2.99%  '

# Row 26: Monero
$ws.Range('D26').Value = '160.84'
$ws.Range('E26').Value = '  +0.74%  '

# Row 27: EthereumClassic
$ws.Range('D27').Value = '20.58'
$ws.Range('E27').Value = '  -2.56%  '

# Row 28: WrappedliquidstakedEther2.0
$ws.Range('D28').Value = '2.010.35'
$ws.Range('E28').Value = '  -1.10%  '

# Row 29: LidoDAOToken
$ws.Range('D29').Value = '2.405'
$ws.Range('E29').Value = '  +0.28%  '

# Row 30: BitcoinCash
$ws.Range('D30').Value = '127.79'
$ws.Range('E30').Value = '  +1.68%  '

# Row 31: Stellar
$ws.Range('D31').Value = '0.1087'
$ws.Range('E31').Value = '  -0.49%  '

# Row 32: ImmutableX
$ws.Range('E32').Value = '  -4.54%  '

# Row 33: HuobiToken
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '3.665'
$ws.Range('E33').Value = '  -0.30%  '

# Row 34: Filecoin
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.584'
$ws.Range('E34').Value = '  -2.59%  '

# Row 35: Hedera
$ws.Range('D35').Value = '0.07166'
$ws.Range('E35').Value = '  -6.79%  '

# Row 36: FraxShare
$ws.Range('D36').Value = '9.116'
$ws.Range('E36').Value = '  +4.13%  '

# Row 37: VeChain
$ws.Range('E37').Value = '  -2.09%  '

# Row 38: Algorand
$ws.Range('D38').Value = '0.2156'
$ws.Range('E38').Value = '  -3.27%  '

# Row 39: Aptos
$ws.Range('D39').Value = '11.66'
$ws.Range('E39').Value = '  +1.32%  '

# Row 40: InternetComputer(DFINITY)
$ws.Range('D40').Value = '5.060'
$ws.Range('E40').Value = '  -4.22%  '

# Row 41: TheSandbox
$ws.Range('D41').Value = '0.6212'
$ws.Range('E41').Value = '  -3.02%  '

# Row 42: Frax
$ws.Range('D42').Value = '0.9973'
$ws.Range('E42').Value = '  -0.51%  '

# Row 43: TrustWalletToken
$ws.Range('E43').Value = '  -3.59%  '

# Row 44: WEMIXTOKEN
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '1.324'
$ws.Range('E44').Value = '  -5.39%  '

# Row 45: EnergySwap
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '13.23'
$ws.Range('E45').Value = '  -2.93%  '

# Row 46: Decentraland
$ws.Range('D46').Value = '0.5976'
$ws.Range('E46').Value = '  -2.89%  '

# Row 47: PancakeSwap
$ws.Range('D47').Value = '3.748'

# Row 48: Quant
$ws.Range('D48').Value = '125.56'
$ws.Range('E48').Value = '  -1.56%  '

# Row 49: EOS
$ws.Range('D49').Value = '1.211'
$ws.Range('E49').Value = '  +0.67%  '

# Row 50: NEARProtocol
$ws.Range('D50').Value = '1.941'
$ws.Range('E50').Value = '  -2.68%  '

# Row 51: Cronos
$ws.Range('D51').Value = '0.06866'
$ws.Range('E51').Value = '  -1.77%  '
